$wb = $excel.ActiveWorkbook

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1023.1667
$ws.Cells.Item(28, 9).Value = 496.84616
$ws.Cells.Item(28, 10).Value = 2391.6
$ws.Cells.Item(28, 11).Value = 496.84616
$ws.Cells.Item(28, 12).Value = 2391.6
$ws.Cells.Item(28, 13).Value = -11.84616
$ws.Cells.Item(28, 14).Value = -3361.6

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 189.78572
$ws.Cells.Item(41, 9).Value = 193.33333
$ws.Cells.Item(41, 11).Value = 193.33333
$ws.Cells.Item(41, 13).Value = 246.66667

# Sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 985.94446
$ws.Cells.Item(53, 9).Value = 1076.4667
$ws.Cells.Item(53, 11).Value = 1076.4667
$ws.Cells.Item(53, 13).Value = -439.4666999999999

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4791.2666
$ws.Cells.Item(116, 9).Value = 4922.1
$ws.Cells.Item(116, 11).Value = 4922.1
$ws.Cells.Item(116, 13).Value = -1480.1

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 9333.764999999999
$ws.Cells.Item(132, 9).Value = 9292.125
$ws.Cells.Item(132, 11).Value = 27876.375
$ws.Cells.Item(132, 13).Value = -25346.375

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2475.9375
$ws.Cells.Item(141, 9).Value = 2441.1333
$ws.Cells.Item(141, 11).Value = 7323.3999
$ws.Cells.Item(141, 13).Value = -2143.3999

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2101.75
$ws.Cells.Item(61, 9).Value = 1602.6364
$ws.Cells.Item(61, 11).Value = 1602.6364
$ws.Cells.Item(61, 13).Value = -1390.6364

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 311656.84
$ws.Cells.Item(74, 9).Value = 558372.8
$ws.Cells.Item(74, 10).Value = 3261.875
$ws.Cells.Item(74, 11).Value = 558372.8
$ws.Cells.Item(74, 12).Value = 3261.875
$ws.Cells.Item(74, 13).Value = -557498.8
$ws.Cells.Item(74, 14).Value = -5009.875

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 311656.84
$ws.Cells.Item(77, 9).Value = 558372.8
$ws.Cells.Item(77, 10).Value = 3261.875
$ws.Cells.Item(77, 11).Value = 2791864
$ws.Cells.Item(77, 12).Value = 16309.375
$ws.Cells.Item(77, 13).Value = -2787496
$ws.Cells.Item(77, 14).Value = -25045.375

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2551.348
$ws.Cells.Item(102, 9).Value = 2174.05
$ws.Cells.Item(102, 11).Value = 2174.05
$ws.Cells.Item(102, 13).Value = -552.0500000000002

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2394.121
$ws.Cells.Item(132, 9).Value = 2127.3333
$ws.Cells.Item(132, 11).Value = 6381.999899999999
$ws.Cells.Item(132, 13).Value = -3851.999899999999

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 2101.75
$ws.Cells.Item(136, 9).Value = 1602.6364
$ws.Cells.Item(136, 11).Value = 4807.9092
$ws.Cells.Item(136, 13).Value = -2257.9092

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 22734986
$ws.Cells.Item(20, 9).Value = 27786486
$ws.Cells.Item(20, 11).Value = 27786486
$ws.Cells.Item(20, 13).Value = -27786239

# Sheet BSM, row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 279.6842
$ws.Cells.Item(80, 9).Value = 203.25
$ws.Cells.Item(80, 10).Value = 300.06668
$ws.Cells.Item(80, 11).Value = 203.25
$ws.Cells.Item(80, 12).Value = 300.06668
$ws.Cells.Item(80, 13).Value = 794.75
$ws.Cells.Item(80, 14).Value = -2296.06668

# Sheet BSM, row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 279.6842
$ws.Cells.Item(83, 9).Value = 203.25
$ws.Cells.Item(83, 10).Value = 300.06668
$ws.Cells.Item(83, 11).Value = 1016.25
$ws.Cells.Item(83, 12).Value = 1500.3334
$ws.Cells.Item(83, 13).Value = 3975.75
$ws.Cells.Item(83, 14).Value = -11484.3334

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2122852.8
$ws.Cells.Item(31, 9).Value = 3439.6572
$ws.Cells.Item(31, 10).Value = 5213663.5
$ws.Cells.Item(31, 11).Value = 3439.6572
$ws.Cells.Item(31, 12).Value = 5213663.5
$ws.Cells.Item(31, 13).Value = -3144.6572
$ws.Cells.Item(31, 14).Value = -5214253.5

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2122852.8
$ws.Cells.Item(34, 9).Value = 3439.6572
$ws.Cells.Item(34, 10).Value = 5213663.5
$ws.Cells.Item(34, 11).Value = 3439.6572
$ws.Cells.Item(34, 12).Value = 5213663.5
$ws.Cells.Item(34, 13).Value = -3237.6572
$ws.Cells.Item(34, 14).Value = -5214067.5

# Sheet CRP, row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()

# Sheet CRP, row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 22229236
$ws.Cells.Item(132, 9).Value = 4911.778
$ws.Cells.Item(132, 11).Value = 14735.334
$ws.Cells.Item(132, 13).Value = -12205.334

# Sheet CRP, row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 69990
$ws.Cells.Item(138, 10).Value = 69990
$ws.Cells.Item(138, 12).Value = 69990
$ws.Cells.Item(138, 14).Value = -80270

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 873.1667
$ws.Cells.Item(5, 9).Value = 873.1667
$ws.Cells.Item(5, 11).Value = 2619.5001
$ws.Cells.Item(5, 13).Value = -2507.5001

# Sheet CUL, row 52
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 537.6667
$ws.Cells.Item(52, 10).Value = 537.6667
$ws.Cells.Item(52, 12).Value = 1613.0001
$ws.Cells.Item(52, 14).Value = -2145.0001

# Sheet CUL, row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 5948
$ws.Cells.Item(56, 9).Value = 5948
$ws.Cells.Item(56, 11).Value = 5948
$ws.Cells.Item(56, 13).Value = -5418

# Sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 2167.875
$ws.Cells.Item(86, 10).Value = 2248.1667
$ws.Cells.Item(86, 12).Value = 6744.500100000001
$ws.Cells.Item(86, 14).Value = -9116.500100000001

# Sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 2167.875
$ws.Cells.Item(89, 10).Value = 2248.1667
$ws.Cells.Item(89, 12).Value = 20233.5003
$ws.Cells.Item(89, 14).Value = -32089.5003

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 873.1667
$ws.Cells.Item(135, 9).Value = 873.1667
$ws.Cells.Item(135, 11).Value = 7858.5003
$ws.Cells.Item(135, 13).Value = -5323.5003

# Sheet GSM, row 49
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 111115950
$ws.Cells.Item(80, 9).Value = 333337340
$ws.Cells.Item(80, 10).Value = 5266.6665
$ws.Cells.Item(80, 11).Value = 333337340
$ws.Cells.Item(80, 12).Value = 5266.6665
$ws.Cells.Item(80, 13).Value = -333336342
$ws.Cells.Item(80, 14).Value = -7262.6665

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 111115950
$ws.Cells.Item(83, 9).Value = 333337340
$ws.Cells.Item(83, 10).Value = 5266.6665
$ws.Cells.Item(83, 11).Value = 1666686700
$ws.Cells.Item(83, 12).Value = 26333.3325
$ws.Cells.Item(83, 13).Value = -1666681708
$ws.Cells.Item(83, 14).Value = -36317.3325

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1829.8
$ws.Cells.Item(102, 9).Value = 1199.3334
$ws.Cells.Item(102, 11).Value = 1199.3334
$ws.Cells.Item(102, 13).Value = 422.6666

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2953.0557
$ws.Cells.Item(132, 9).Value = 2892.9092
$ws.Cells.Item(132, 11).Value = 8678.7276
$ws.Cells.Item(132, 13).Value = -6148.7276

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2171.8
$ws.Cells.Item(7, 9).Value = 2091.7334
$ws.Cells.Item(7, 10).Value = 2412
$ws.Cells.Item(7, 11).Value = 2091.7334
$ws.Cells.Item(7, 12).Value = 2412
$ws.Cells.Item(7, 13).Value = -1979.7334
$ws.Cells.Item(7, 14).Value = -2636

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3187.5789
$ws.Cells.Item(46, 9).Value = 2524.0908
$ws.Cells.Item(46, 10).Value = 4099.875
$ws.Cells.Item(46, 11).Value = 2524.0908
$ws.Cells.Item(46, 12).Value = 4099.875
$ws.Cells.Item(46, 13).Value = -2336.0908
$ws.Cells.Item(46, 14).Value = -4475.875

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3086.1428
$ws.Cells.Item(68, 9).Value = 2964.6
$ws.Cells.Item(68, 10).Value = 3390
$ws.Cells.Item(68, 11).Value = 2964.6
$ws.Cells.Item(68, 12).Value = 3390
$ws.Cells.Item(68, 13).Value = -2215.6
$ws.Cells.Item(68, 14).Value = -4888

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 3086.1428
$ws.Cells.Item(71, 9).Value = 2964.6
$ws.Cells.Item(71, 10).Value = 3390
$ws.Cells.Item(71, 11).Value = 14823
$ws.Cells.Item(71, 12).Value = 16950
$ws.Cells.Item(71, 13).Value = -11079
$ws.Cells.Item(71, 14).Value = -24438

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1554.4
$ws.Cells.Item(93, 9).Value = 1554.4
$ws.Cells.Item(93, 11).Value = 1554.4
$ws.Cells.Item(93, 13).Value = -306.4000000000001

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 10000
$ws.Cells.Item(100, 9).Value = 10000
$ws.Cells.Item(100, 11).Value = 10000
$ws.Cells.Item(100, 13).Value = -9459

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 2171.8
$ws.Cells.Item(126, 9).Value = 2091.7334
$ws.Cells.Item(126, 10).Value = 2412
$ws.Cells.Item(126, 11).Value = 6275.2002
$ws.Cells.Item(126, 12).Value = 7236
$ws.Cells.Item(126, 13).Value = -3805.2002
$ws.Cells.Item(126, 14).Value = -12176

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6099.769
$ws.Cells.Item(132, 9).Value = 2225.875
$ws.Cells.Item(132, 10).Value = 12298
$ws.Cells.Item(132, 11).Value = 6677.625
$ws.Cells.Item(132, 12).Value = 36894
$ws.Cells.Item(132, 13).Value = -4147.625
$ws.Cells.Item(132, 14).Value = -41954

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5213.4
$ws.Cells.Item(136, 9).Value = 1651.2667
$ws.Cells.Item(136, 10).Value = 15899.8
$ws.Cells.Item(136, 11).Value = 4953.800099999999
$ws.Cells.Item(136, 12).Value = 47699.39999999999
$ws.Cells.Item(136, 13).Value = -2403.800099999999
$ws.Cells.Item(136, 14).Value = -52799.39999999999

# Sheet WVR, row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 19521
$ws.Cells.Item(52, 9).Value = 19521
$ws.Cells.Item(52, 11).Value = 19521
$ws.Cells.Item(52, 13).Value = -19295

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 792.0357
$ws.Cells.Item(113, 9).Value = 755.9524
$ws.Cells.Item(113, 11).Value = 2267.8572
$ws.Cells.Item(113, 13).Value = -97.85719999999992
